$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new PatientID values to the first empty rows below the existing data
$ws.Range("A9").Value = "PEP_ID-2009237"
$ws.Range("A10").Value = "PEP_ID-2009241"
$ws.Range("A11").Value = "PEP_ID-2009243"
